$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-shuffles the Fecha/Volumen/Precio/Origen/Precio-Kg data block
# (columns D, J, K, L, M, O, P) across the data rows (2-25), leaving the
# Mercado/Region/Categoria/etc. columns and row 11 untouched.

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the rows that move, keyed by source row.
$snapshot = @{}
foreach ($r in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)) {
    $rowData = @{}
    foreach ($c in $cols) {
        $addr = "$c$r"
        $rowData[$c] = $ws.Range($addr).Value2
    }
    $snapshot[$r] = $rowData
}

# Write snapshot values into their new destination rows.
$destMap = @{
    2 = 24
    3 = 16
    4 = 21
    5 = 25
    6 = 7
    7 = 9
    8 = 13
    9 = 20
    10 = 4
    12 = 19
    13 = 3
    14 = 22
    15 = 10
    16 = 12
    17 = 2
    18 = 6
    19 = 15
    20 = 23
    21 = 8
    22 = 17
    23 = 18
    24 = 5
    25 = 14
}

foreach ($dest in $destMap.Keys) {
    $src = $destMap[$dest]
    $rowData = $snapshot[$src]
    foreach ($c in $cols) {
        $addr = "$c$dest"
        $ws.Range($addr).Value2 = $rowData[$c]
    }
}